$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell F1 - copy the style of E1 (bold header style) then set the text
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy() | Out-Null
$ws.Range("F1").PasteSpecial(-4122) | Out-Null

# time_taken values for data rows 2-42
$timeTaken = @(
    "2021-10-05 13:42:08.134505",
    "2021-10-05 13:42:08.134517",
    "2021-10-05 13:42:08.134520",
    "2021-10-05 13:42:08.134523",
    "2021-10-05 13:42:08.134526",
    "2021-10-05 13:42:08.134529",
    "2021-10-05 13:42:08.134532",
    "2021-10-05 13:42:08.134534",
    "2021-10-05 13:42:08.134537",
    "2021-10-05 13:42:08.134540",
    "2021-10-05 13:42:08.134543",
    "2021-10-05 13:42:08.134545",
    "2021-10-05 13:42:08.134548",
    "2021-10-05 13:42:08.134550",
    "2021-10-05 13:42:08.134553",
    "2021-10-05 13:42:08.134556",
    "2021-10-05 13:42:08.134559",
    "2021-10-05 13:42:08.134562",
    "2021-10-05 13:42:08.134565",
    "2021-10-05 13:42:08.134568",
    "2021-10-05 13:42:08.134571",
    "2021-10-05 13:42:08.134573",
    "2021-10-05 13:42:08.134576",
    "2021-10-05 13:42:08.134578",
    "2021-10-05 13:42:08.134581",
    "2021-10-05 13:42:08.134584",
    "2021-10-05 13:42:08.134586",
    "2021-10-05 13:42:08.134589",
    "2021-10-05 13:42:08.134591",
    "2021-10-05 13:42:08.134594",
    "2021-10-05 13:42:08.134597",
    "2021-10-05 13:42:08.134599",
    "2021-10-05 13:42:08.134602",
    "2021-10-05 13:42:08.134605",
    "2021-10-05 13:42:08.134608",
    "2021-10-05 13:42:08.134610",
    "2021-10-05 13:42:08.134613",
    "2021-10-05 13:42:08.134616",
    "2021-10-05 13:42:08.134618",
    "2021-10-05 13:42:08.134621",
    "2021-10-05 13:42:08.134624"
)

for ($i = 0; $i -lt $timeTaken.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timeTaken[$i]
}
